$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Replace the old value with the new one inside a single table cell,
# addressed by (row, column). The cell's Range.Start/End is recomputed
# right before the Find/Replace call (not cached across edits) because
# re-reading Table.Cell() positions after earlier mutations — or batching
# several captured ranges before applying any of them — proved unreliable
# whenever two cells happen to share identical starting text (this sheet
# has "54÷8=6, 6" twice). Doing it immediately-before-use, one cell at a
# time, keeps every offset accurate for the cell actually being edited.
function Set-CellText($row, $col, $old, $new) {
    $cell = $t.Cell($row, $col)
    $r = $d.Range($cell.Range.Start, $cell.Range.End)
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Update the date heading above the table.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Find.Execute("2023-12-03 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-04 Monday", 2) | Out-Null

# Update every division-problem cell, row by row, left to right.
Set-CellText 1 1 "65÷7=9, 2" "23÷5=4, 3"
Set-CellText 1 2 "45÷4=11, 1" "81÷6=13, 3"
Set-CellText 1 3 "46÷8=5, 6" "28÷2=14, 0"
Set-CellText 1 4 "10÷8=1, 2" "69÷6=11, 3"
Set-CellText 1 5 "16÷4=4, 0" "74÷9=8, 2"

Set-CellText 5 1 "41÷9=4, 5" "64÷4=16, 0"
Set-CellText 5 2 "58÷8=7, 2" "59÷6=9, 5"
Set-CellText 5 3 "43÷2=21, 1" "46÷7=6, 4"
Set-CellText 5 4 "54÷8=6, 6" "11÷9=1, 2"
Set-CellText 5 5 "13÷4=3, 1" "86÷3=28, 2"

Set-CellText 9 1 "53÷7=7, 4" "45÷4=11, 1"
Set-CellText 9 2 "15÷7=2, 1" "46÷4=11, 2"
Set-CellText 9 3 "91÷7=13, 0" "37÷8=4, 5"
Set-CellText 9 4 "39÷2=19, 1" "42÷2=21, 0"
Set-CellText 9 5 "91÷8=11, 3" "98÷5=19, 3"

Set-CellText 13 1 "75÷9=8, 3" "85÷3=28, 1"
Set-CellText 13 2 "60÷3=20, 0" "84÷9=9, 3"
Set-CellText 13 3 "60÷9=6, 6" "68÷5=13, 3"
Set-CellText 13 4 "86÷5=17, 1" "29÷4=7, 1"
Set-CellText 13 5 "54÷8=6, 6" "22÷3=7, 1"

Set-CellText 17 1 "14÷4=3, 2" "13÷7=1, 6"
Set-CellText 17 2 "88÷9=9, 7" "50÷6=8, 2"
Set-CellText 17 3 "90÷6=15, 0" "81÷7=11, 4"
Set-CellText 17 4 "73÷5=14, 3" "41÷6=6, 5"
Set-CellText 17 5 "15÷2=7, 1" "53÷4=13, 1"

Write-Output "done"
